$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Octubre de 2020 a las 06:12"

# 2) Row 5 - India: refresh counts
$ws.Range("B5").Value = 7053806
$ws.Range("C5").Value = 2263
$ws.Range("D5").Value = 6077976
$ws.Range("E5").Value = 867459

# 3) Row 26 - Pakistan: refresh counts
$ws.Range("B26").Value = 318932
$ws.Range("C26").Value = 666
$ws.Range("D26").Value = 303458
$ws.Range("E26").Value = 8904
$ws.Range("G26").Value = 12
$ws.Range("H26").Value = 6570

# 4) Rows 31-33 - Belgica jumps above Rumania & Marruecos with refreshed counts;
#    Rumania and Marruecos keep their previous figures but shift down a row each.
$ws.Range("A31").Value = "Belgica"
$ws.Range("B31").Value = 156931
$ws.Range("C31").Value = 7950
$ws.Range("D31").Value = 20202
$ws.Range("E31").Value = 126554
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 24
$ws.Range("H31").Value = 10175

$ws.Range("A32").Value = "Rumania"
$ws.Range("B32").Value = 152403
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 116628
$ws.Range("E32").Value = 30417
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 5358

$ws.Range("A33").Value = "Marruecos"
$ws.Range("B33").Value = 149841
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 124854
$ws.Range("E33").Value = 22415
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 2572

# 5) Row 54 - Honduras: refresh counts
$ws.Range("B54").Value = 83146
$ws.Range("C54").Value = 594
$ws.Range("D54").Value = 31931
$ws.Range("E54").Value = 48711
$ws.Range("G54").Value = 12
$ws.Range("H54").Value = 2504

# 6) Row 143 - Tailandia: refresh counts
$ws.Range("B143").Value = 3636
$ws.Range("C143").Value = 2
$ws.Range("D143").Value = 3451
$ws.Range("E143").Value = 126

# 7) Row 172 - San Martin (Parte Holandesa): refresh counts
$ws.Range("B172").Value = 707
$ws.Range("C172").Value = 4
$ws.Range("E172").Value = 63

# 8) Row 187 - Butan: refresh counts
$ws.Range("D187").Value = 287
$ws.Range("E187").Value = 19
